$d = $word.ActiveDocument

$pairs = @(
    @("39×86=3354", "71×43=3053"),
    @("77×75=5775", "31×48=1488"),
    @("55×20=1100", "30×33=990"),
    @("12×11=132", "89×59=5251"),
    @("61×24=1464", "82×55=4510"),
    @("31×52=1612", "91×58=5278"),
    @("24×16=384", "70×32=2240"),
    @("54×63=3402", "43×13=559"),
    @("11×24=264", "44×99=4356"),
    @("91×38=3458", "71×25=1775"),
    @("39×29=1131", "90×71=6390"),
    @("77×82=6314", "27×46=1242"),
    @("74×38=2812", "31×46=1426"),
    @("29×30=870", "65×69=4485"),
    @("89×23=2047", "79×64=5056"),
    @("85×83=7055", "53×62=3286"),
    @("34×58=1972", "90×26=2340"),
    @("93×30=2790", "15×18=270"),
    @("88×28=2464", "13×23=299"),
    @("62×57=3534", "73×88=6424"),
    @("47×46=2162", "65×72=4680"),
    @("15×33=495", "12×17=204"),
    @("29×57=1653", "59×40=2360"),
    @("81×88=7128", "56×28=1568"),
    @("89×56=4984", "52×58=3016")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
